# Regenerate save_data column G ("K") values: use K instead of Strike#,
# regen std/mean, calc and write s_vals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 2
    6  = 1
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 2
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 1
    18 = 3
    19 = 1
    20 = 2
    21 = 1
    22 = 1
    24 = 1
    25 = 0
    26 = 2
    27 = 0
    28 = 1
    29 = 1
    30 = 2
    31 = 1
    32 = 2
    33 = 2
    34 = 1
    35 = 0
    36 = 1
    37 = 0
    38 = 2
    39 = 1
    40 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
